$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: headers (Indonesian -> English translation) ---
$ws.Range("A1").Value = "Component"
# B1 "Product ID" stays unchanged
$ws.Range("C1").Value = "Amount"

# --- Row 2: Service 102 (was Jasa 101 / 101) ---
$ws.Range("A2").Value = "Service 102"
$ws.Range("B2").Value = "'102"
$ws.Range("B2").Style = "Normal"
$ws.Range("C2").Value = 218.88

# --- Row 3: Service 103 (was Jasa 102 / 102) ---
$ws.Range("A3").Value = "Service 103"
$ws.Range("B3").Value = "'103"
$ws.Range("B3").Style = "Normal"
$ws.Range("C3").Value = 713.6999999999999

# --- Row 4: Service 101 (was Jasa 103 / 103) ---
$ws.Range("A4").Value = "Service 101"
$ws.Range("B4").Value = "'101"
$ws.Range("B4").Style = "Normal"
$ws.Range("C4").Value = 2144.84

# --- Row 5: Total Sales (was Total Pendapatan) ---
$ws.Range("A5").Value = "Total Sales"
# C5 value (50967.56) unchanged

# --- Row 6: Daily Profit and Loss (was Laba Rugi Harian) ---
$ws.Range("A6").Value = "Daily Profit and Loss"
# B6 value unchanged

# --- Row 7: Total Cash Flow (was Total Arus Kas) ---
$ws.Range("A7").Value = "Total Cash Flow"
# B7 value unchanged

# --- Row 8: Daily Cash Balance (was Saldo Kas Harian) ---
$ws.Range("A8").Value = "Daily Cash Balance"
# B8 value unchanged

# --- Row 9: Total Assets (was Total Aset) ---
$ws.Range("A9").Value = "Total Assets"
$ws.Range("B9").Value = 10093.75999999999

# --- Row 10: Total Debt (was Total Utang) ---
$ws.Range("A10").Value = "Total Debt"
$ws.Range("B10").Value = 800000

# --- Row 11: Total Equity (was Total Ekuitas) ---
$ws.Range("A11").Value = "Total Equity"
$ws.Range("B11").Value = -789906.24
